$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for cells whose new values look like plain numbers,
# so they stay text (matching the original inlineStr cells) instead of
# being auto-converted to numeric values by Excel.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D32", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin values (price + 1h volume change).
$ws.Range("D2").Value = "69.935.40"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "3.796.16"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "621.28"
$ws.Range("E5").Value = "  +3.91%  "
$ws.Range("D6").Value = "176.97"
$ws.Range("E6").Value = "  -3.81%  "
$ws.Range("D7").Value = "3.792.49"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "6.28"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").Value = "0.492"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "40.68"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "0.0000261"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "4.425.41"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "3.791.35"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "69.986.75"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "16.80"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "510.02"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "9.55"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "87.67"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").Value = "13.21"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0000141"
$ws.Range("E27").Value = "  +29.04%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.97"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("D37").Value = "6.21"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "0.134"
$ws.Range("E38").Value = "  +5.57%  "
$ws.Range("D39").Value = "0.331"
$ws.Range("E39").Value = "  -3.12%  "
$ws.Range("D40").Value = "2.13"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "50.94"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "45.18"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").Value = "8.72"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "415.81"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "3.034.18"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "0.0362"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "27.32"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D50").Value = "137.94"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("D51").Value = "2.46"
$ws.Range("E51").Value = "  +0.98%  "

Write-Output "Applied 94 cell updates"
